$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (index 1)
# Rows were re-ordered/re-generated by a fresh report run:
#   old row2 (3a6762b9...) content moved down to row4 (status updated)
#   old row3 (ffff8c1cf609...) content moved up to row2
#   old row4 (ffffffdf1752a0...) content moved up to row3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

$wsOverview.Range("A2").Value = "ffff8c1cf609-4f2a-4bde-928b-57d98bb3639e.md"
$wsOverview.Range("B2").Value = "e2e\ffff8c1cf609-4f2a-4bde-928b-57d98bb3639e.md"
$wsOverview.Range("G2").Value = "2016-08-25 15:06:47"

$wsOverview.Range("A3").Value = "ffffffdf1752a0-5523-4ffe-b5e6-4fd6181af3b0.md"
$wsOverview.Range("B3").Value = "e2e\ffffffdf1752a0-5523-4ffe-b5e6-4fd6181af3b0.md"

$wsOverview.Range("A4").Value = "3a6762b9-1712-4999-9677-c964c0d67906.md"
$wsOverview.Range("B4").Value = "e2e\3a6762b9-1712-4999-9677-c964c0d67906.md"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-25 15:09:25"

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = "e2e\ffff8c1cf609-4f2a-4bde-928b-57d98bb3639e.md"
    } elseif ($addr -eq '$B$3') {
        $h.TextToDisplay = "e2e\ffffffdf1752a0-5523-4ffe-b5e6-4fd6181af3b0.md"
    } elseif ($addr -eq '$B$4') {
        $h.TextToDisplay = "e2e\3a6762b9-1712-4999-9677-c964c0d67906.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (index 2)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item(2)

$wsZh.Range("A2").Value = "ffff8c1cf609-4f2a-4bde-928b-57d98bb3639e.md"
$wsZh.Range("G2").Value = "5fc316cb-1163-4287-adb3-ae991cadfad9.e8ba3ad45d87ef6fc7e4b43cb4485f17c8599e7a.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-25 15:06:42"
$wsZh.Range("I2").Value = "5fc316cb-1163-4287-adb3-ae991cadfad9.md"
$wsZh.Range("J2").Value = "5fc316cb-1163-4287-adb3-ae991cadfad9.e8ba3ad45d87ef6fc7e4b43cb4485f17c8599e7a.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-25 15:07:01"

$wsZh.Range("A3").Value = "ffffffdf1752a0-5523-4ffe-b5e6-4fd6181af3b0.md"
$wsZh.Range("F3").Value = "True"

$wsZh.Range("A4").Value = "3a6762b9-1712-4999-9677-c964c0d67906.md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = "3a6762b9-1712-4999-9677-c964c0d67906.d116c5551b9ae572883359643007996b07dcb089.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-08-25 15:09:20"
$wsZh.Range("I4").Value = "3a6762b9-1712-4999-9677-c964c0d67906.md"
$wsZh.Range("J4").Value = "3a6762b9-1712-4999-9677-c964c0d67906.d116c5551b9ae572883359643007996b07dcb089.zh-cn.xlf"
$wsZh.Range("K4").Value = "2016-08-25 15:08:42"
$wsZh.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6dbe1b027c3e3a6f539e8a337267deea476cdbb5/e2e/3a6762b9-1712-4999-9677-c964c0d67906.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7abd55954c738ddce51287c85d7a5d287297f923/e2e/3a6762b9-1712-4999-9677-c964c0d67906.md."

foreach ($h in $wsZh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "ffff8c1cf609-4f2a-4bde-928b-57d98bb3639e.md"
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = "5fc316cb-1163-4287-adb3-ae991cadfad9.md"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "ffffffdf1752a0-5523-4ffe-b5e6-4fd6181af3b0.md"
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = "5fc316cb-1163-4287-adb3-ae991cadfad9.md"
    } elseif ($addr -eq '$A$4') {
        $h.TextToDisplay = "3a6762b9-1712-4999-9677-c964c0d67906.md"
    } elseif ($addr -eq '$I$4') {
        $h.TextToDisplay = "3a6762b9-1712-4999-9677-c964c0d67906.md"
    }
}

# Error Detail column widened to fit the new long diagnostic text.
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# Sheet "de-de" (index 3)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item(3)

$wsDe.Range("A2").Value = "ffff8c1cf609-4f2a-4bde-928b-57d98bb3639e.md"
$wsDe.Range("G2").Value = "5fc316cb-1163-4287-adb3-ae991cadfad9.e8ba3ad45d87ef6fc7e4b43cb4485f17c8599e7a.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-25 15:06:47"
$wsDe.Range("I2").Value = "5fc316cb-1163-4287-adb3-ae991cadfad9.md"
$wsDe.Range("J2").Value = "5fc316cb-1163-4287-adb3-ae991cadfad9.e8ba3ad45d87ef6fc7e4b43cb4485f17c8599e7a.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-25 15:07:18"

$wsDe.Range("A3").Value = "ffffffdf1752a0-5523-4ffe-b5e6-4fd6181af3b0.md"
$wsDe.Range("F3").Value = "True"

$wsDe.Range("A4").Value = "3a6762b9-1712-4999-9677-c964c0d67906.md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = "3a6762b9-1712-4999-9677-c964c0d67906.d116c5551b9ae572883359643007996b07dcb089.de-de.xlf"
$wsDe.Range("H4").Value = "2016-08-25 15:09:25"
$wsDe.Range("I4").Value = "3a6762b9-1712-4999-9677-c964c0d67906.md"
$wsDe.Range("J4").Value = "3a6762b9-1712-4999-9677-c964c0d67906.d116c5551b9ae572883359643007996b07dcb089.de-de.xlf"
$wsDe.Range("K4").Value = "2016-08-25 15:08:49"
$wsDe.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6dbe1b027c3e3a6f539e8a337267deea476cdbb5/e2e/3a6762b9-1712-4999-9677-c964c0d67906.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7abd55954c738ddce51287c85d7a5d287297f923/e2e/3a6762b9-1712-4999-9677-c964c0d67906.md."

foreach ($h in $wsDe.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "ffff8c1cf609-4f2a-4bde-928b-57d98bb3639e.md"
    } elseif ($addr -eq '$I$2') {
        $h.TextToDisplay = "5fc316cb-1163-4287-adb3-ae991cadfad9.md"
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = "ffffffdf1752a0-5523-4ffe-b5e6-4fd6181af3b0.md"
    } elseif ($addr -eq '$I$3') {
        $h.TextToDisplay = "5fc316cb-1163-4287-adb3-ae991cadfad9.md"
    } elseif ($addr -eq '$A$4') {
        $h.TextToDisplay = "3a6762b9-1712-4999-9677-c964c0d67906.md"
    } elseif ($addr -eq '$I$4') {
        $h.TextToDisplay = "3a6762b9-1712-4999-9677-c964c0d67906.md"
    }
}

# Error Detail column widened to fit the new long diagnostic text.
$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664
